# Applies the 9/17-afternoon Matlab-import update to PS1/q2Matlab.xlsx:
#   - relabel the row-index column header "Variable" -> "Row"
#   - replace the (near-singular) se/t_stat/p_val/CI columns with the
#     corrected "Power calculations" figures, on both the symm and chol sheets
#   - columns C/D/F/G re-fit their width to the new numbers (col E is unaffected
#     since its formatted width does not change)

$wb = $excel.ActiveWorkbook
$wsSymm = $wb.Worksheets.Item("symm")
$wsChol = $wb.Worksheets.Item("chol")

# --- Header rename: "Variable" -> "Row" (A1 on both sheets, shared string) ---
$wsSymm.Range("A1").Value = "Row"
$wsChol.Range("A1").Value = "Row"

# --- Updated regression statistics: "symm" sheet (rows 2-11) ---
$wsSymm.Range("C2").Value = 94136.76726107973
$wsSymm.Range("D2").Value = 1.453341734805656
$wsSymm.Range("E2").Value = 0.07341765273127265
$wsSymm.Range("F2").Value = -178022.51072620126
$wsSymm.Range("G2").Value = 190993.6169372313
$wsSymm.Range("C3").Value = 13311.508560303344
$wsSymm.Range("D3").Value = 2.433310825400328
$wsSymm.Range("E3").Value = 0.007677568339422791
$wsSymm.Range("F3").Value = -24555.074322749573
$wsSymm.Range("G3").Value = 27626.039233639534
$wsSymm.Range("C4").Value = 16581.021845903015
$wsSymm.Range("D4").Value = -3.298122738534785
$wsSymm.Range("E4").Value = 0.9994741991683567
$wsSymm.Range("F4").Value = -35091.17946025448
$wsSymm.Range("G4").Value = 29906.42617568533
$wsSymm.Range("C5").Value = 844.0701183206412
$wsSymm.Range("D5").Value = 0.9831993697442545
$wsSymm.Range("E5").Value = 0.1630223872707217
$wsSymm.Range("F5").Value = -1615.0369104956087
$wsSymm.Range("G5").Value = 1693.7179533213045
$wsSymm.Range("C6").Value = 19702.836261608834
$wsSymm.Range("D6").Value = -0.7928659301812048
$wsSymm.Range("E6").Value = 0.7858603017517667
$wsSymm.Range("F6").Value = -39358.09903417026
$wsSymm.Range("G6").Value = 37877.01911133636
$wsSymm.Range("C7").Value = 1121.428443250967
$wsSymm.Range("D7").Value = 1.1301998184381536
$wsSymm.Range("E7").Value = 0.1295011014247225
$wsSymm.Range("F7").Value = -2137.9174148317593
$wsSymm.Range("G7").Value = 2258.082082712031
$wsSymm.Range("C8").Value = 2.1627850076780044
$wsSymm.Range("D8").Value = -0.2912215907236558
$wsSymm.Range("E8").Value = 0.6144910763936957
$wsSymm.Range("F8").Value = -4.268916353864959
$wsSymm.Range("G8").Value = 4.209200876232818
$wsSymm.Range("C9").Value = 2.749041283811572
$wsSymm.Range("D9").Value = 1.3456003236602125
$wsSymm.Range("E9").Value = 0.089559094438322
$wsSymm.Range("F9").Value = -5.212766245075933
$wsSymm.Range("G9").Value = 5.56347558746543
$wsSymm.Range("C10").Value = 31408.59706078061
$wsSymm.Range("D10").Value = 0.8838893632226543
$wsSymm.Range("E10").Value = 0.18861714753315018
$wsSymm.Range("F10").Value = -60244.818267448165
$wsSymm.Range("G10").Value = 62876.88221081182
$wsSymm.Range("C11").Value = 26600.90169851313
$wsSymm.Range("D11").Value = -0.9259992035496175
$wsSymm.Range("E11").Value = 0.8225252838935417
$wsSymm.Range("F11").Value = -53305.45572925461
$wsSymm.Range("G11").Value = 50970.07892891686

# --- Updated regression statistics: "chol" sheet (rows 2-11) ---
$wsChol.Range("C2").Value = 94136.76726102417
$wsChol.Range("D2").Value = 1.4533417348065871
$wsChol.Range("E2").Value = 0.07341765273114387
$wsChol.Range("F2").Value = -178022.51072609203
$wsChol.Range("G2").Value = 190993.61693712272
$wsChol.Range("C3").Value = 13311.50856030318
$wsChol.Range("D3").Value = 2.433310825400363
$wsChol.Range("E3").Value = 0.007677568339422014
$wsChol.Range("F3").Value = -24555.07432274925
$wsChol.Range("G3").Value = 27626.039233639218
$wsChol.Range("C4").Value = 16581.021845903037
$wsChol.Range("D4").Value = -3.298122738534806
$wsChol.Range("E4").Value = 0.9994741991683568
$wsChol.Range("F4").Value = -35091.17946025454
$wsChol.Range("G4").Value = 29906.426175685356
$wsChol.Range("C5").Value = 844.0701183205313
$wsChol.Range("D5").Value = 0.9831993697443856
$wsChol.Range("E5").Value = 0.1630223872706894
$wsChol.Range("F5").Value = -1615.0369104953934
$wsChol.Range("G5").Value = 1693.7179533210892
$wsChol.Range("C6").Value = 19702.836261597535
$wsChol.Range("D6").Value = -0.792865930181765
$wsChol.Range("E6").Value = 0.7858603017519298
$wsChol.Range("F6").Value = -39358.09903414822
$wsChol.Range("G6").Value = 37877.01911131411
$wsChol.Range("C7").Value = 1121.4284432504242
$wsChol.Range("D7").Value = 1.130199818438733
$wsChol.Range("E7").Value = 0.12950110142460092
$wsChol.Range("F7").Value = -2137.917414830694
$wsChol.Range("G7").Value = 2258.082082710969
$wsChol.Range("C8").Value = 2.162785007678042
$wsChol.Range("D8").Value = -0.29122159072364634
$wsChol.Range("E8").Value = 0.6144910763936923
$wsChol.Range("F8").Value = -4.268916353865033
$wsChol.Range("G8").Value = 4.2092008762328925
$wsChol.Range("C9").Value = 2.74904128381159
$wsChol.Range("D9").Value = 1.3456003236602017
$wsChol.Range("E9").Value = 0.08955909443832333
$wsChol.Range("F9").Value = -5.212766245075969
$wsChol.Range("G9").Value = 5.563475587465463
$wsChol.Range("C10").Value = 31408.597060780816
$wsChol.Range("D10").Value = 0.8838893632226505
$wsChol.Range("E10").Value = 0.1886171475331515
$wsChol.Range("F10").Value = -60244.81826744857
$wsChol.Range("G10").Value = 62876.882210812226
$wsChol.Range("C11").Value = 26600.901698513204
$wsChol.Range("D11").Value = -0.9259992035496062
$wsChol.Range("E11").Value = 0.8225252838935391
$wsChol.Range("F11").Value = -53305.45572925474
$wsChol.Range("G11").Value = 50970.07892891702

# --- Column width re-fit (bestFit-style autosize after the value change) ---
# Target (real-Excel bestFit) widths are 9.140625 / 6 / 10.85546875 / 10.140625
# for columns C/D/F/G; column widths set through this host are quantized to
# 1/6-character increments, so we pick the ColumnWidth input that lands on
# the closest achievable increment to the real target for each column.
$wsSymm.Columns.Item(3).ColumnWidth = 8.33
$wsSymm.Columns.Item(4).ColumnWidth = 5.165
$wsSymm.Columns.Item(6).ColumnWidth = 10.0
$wsSymm.Columns.Item(7).ColumnWidth = 9.33

$wsChol.Columns.Item(3).ColumnWidth = 8.33
$wsChol.Columns.Item(4).ColumnWidth = 5.165
$wsChol.Columns.Item(6).ColumnWidth = 10.0
$wsChol.Columns.Item(7).ColumnWidth = 9.33
